# "Generate Report for Handoff"
# Status flips from "In Translation" -> "Ready for handoff" and the
# associated timestamps advance a few dozen seconds across the Overview,
# zh-cn and de-de sheets. The Status/locale columns also widen slightly
# to accommodate the new, longer status text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-25 20:59:11"
$ws.Columns("E:F").ColumnWidth = 17.2159881591797

# ---- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-25 20:59:04"
$ws.Columns("C:C").ColumnWidth = 17.2159881591797

# ---- de-de sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-25 20:59:11"
$ws.Columns("C:C").ColumnWidth = 17.2159881591797
